# Fixed a bug in symbolTrigger
# Reorders the data rows (2-21) on the active worksheet to the corrected
# symbol ordering, while leaving headers (row 1) and the trailing
# aggregate rows (22-26) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(301,  6, 45, 30, 60, 45),
    @(1202, 2, 10, 10, 10, 10),
    @(902,  1,  0,  0,  0,  0),
    @(501,  9, 52, 30, 75, 45),
    @(401,  9, 48, 67, 75, 45),
    @(201,  9, 30, 15, 45, 30),
    @(1203, 3, 15, 15, 15, 15),
    @(901, 16, 15, 45, 60, 60),
    @(701,  3, 90, 45, 97, 15),
    @(601,  9, 60, 67, 60, 42),
    @(801,  3, 67, 65, 52, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(101,  9, 30, 15, 60, 15),
    @(1001,18, 30, 75, 60, 72),
    @(2,    0,  2,  2,  2,  2),
    @(1,    0,  2,  2,  2,  2),
    @(502,  0,  4,  0,  0,  0),
    @(1101, 0, 15, 30, 30,  0),
    @(3,    0,  3,  3,  3,  3),
    @(802,  0,  4,  5,  4,  0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 0; $col -lt $values.Length; $col++) {
        $ws.Cells.Item($row, $col + 1).Value = $values[$col]
    }
}
